$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.149822592735291
$ws.Range("B1").Value = 1.042942404747009
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 1.995120406150818
$ws.Range("E1").Value = 0.9888195395469666
